$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 value (OTP column) - login flow now uses a generated code instead of ABC123
$ws.Range("B2").Value = "X1Y2Z3"

# Update the active selection to B2
$ws.Range("B2").Select()
